$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 73, shifting the existing rows 73:89 down to 74:90.
$ws.Rows.Item(73).Insert()

# Populate the newly inserted row 73 with the new weekly price record.
$ws.Cells.Item(73, 1).Value = 4
$ws.Cells.Item(73, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(73, 3).Value = "Los Lagos"
$ws.Cells.Item(73, 4).Value = 45204
$ws.Cells.Item(73, 5).Value = 10
$ws.Cells.Item(73, 6).Value = 100112012
$ws.Cells.Item(73, 7).Value = "Espinaca"
$ws.Cells.Item(73, 8).Value = "Sin especificar"
$ws.Cells.Item(73, 9).Value = "Primera"
$ws.Cells.Item(73, 10).Value = 25
$ws.Cells.Item(73, 11).Value = 13000
$ws.Cells.Item(73, 12).Value = 13000
$ws.Cells.Item(73, 13).Value = 13000
$ws.Cells.Item(73, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(73, 15).Value = "Región Metropolitana"
$ws.Cells.Item(73, 16).Value = 1300
$ws.Cells.Item(73, 17).Value = 10
$ws.Cells.Item(73, 18).Value = "Hortaliza"
